$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 / J1 - copy formatting from existing header cell H1
# then set their text (xlPasteFormats = -4122)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-37: I = 1 (constant), J = same value as H (copy of IP)
for ($r = 2; $r -le 37; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
